$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Pre Consolidation Identifiers -> state id replaced by state fingerprint id)
$ws.Range("A5").Value = "Person State Fingerprint ID"
$ws.Range("B5").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C5").Value = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PreConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"

# Row 6 mapping prefix rename only (CHcr-doc/CHcr-ext -> chc-report-doc/chc-report-ext)
$ws.Range("C6").Value = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PreConsolidationIdentifiers/j:PersonFBIIdentification/nc:IdentificationID"

# Row 8 (Post Consolidation Identifiers -> state id replaced by state fingerprint id)
$ws.Range("A8").Value = "Person State Fingerprint ID"
$ws.Range("B8").Value = "An identification of a person based on a Fingerprint ID."
$ws.Range("C8").Value = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"

# Row 9 takes over the old row 10 (FBI identification) content/format
$ws.Range("A9").Value = "Person FBI Identification ID"
$ws.Range("B9").Value = "A number issued by the FBI's Automated Fingerprint Identification System (AFIS) based on submitted fingerprints."
$ws.Range("C9").Value = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PostConsolidationIdentifiers/j:PersonFBIIdentification/nc:IdentificationID"
$ws.Rows.Item(9).RowHeight = 28

# Old row 10 (duplicate FBI row) is removed entirely
$ws.Rows.Item(10).Delete()

# Update the active selection to match the post-edit state
$ws.Range("C9").Select()
